$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 44576.55555555555
$ws.Range("F2").Value = 44577

$ws.Range("E3").Value = 44576.55555555555
$ws.Range("F3").Value = 44577

$ws.Range("E4").Value = 44576.55555555555
$ws.Range("F4").Value = 44577

$ws.Range("E5").Value = 44576.55555555555
$ws.Range("F5").Value = 44577

$ws.Range("E6").Value = 44576.55555555555
$ws.Range("F6").Value = 44577

$ws.Range("E7").Value = 44576.55555555555
$ws.Range("F7").Value = 44577

$ws.Range("E8").Value = 44571.47986111111
$ws.Range("F8").Value = 44571.76597222222

$ws.Range("E9").Value = 44573.97638888889
$ws.Range("F9").Value = 44574.2625

$ws.Range("B10").Value = "139"
$ws.Range("E10").Value = 44574.28333333333
$ws.Range("F10").Value = 44574.56736111111

$ws.Range("E11").Value = 44574.84722222222
$ws.Range("F11").Value = 44575.16666666666

$ws.Range("E12").Value = 44574.57569444444
$ws.Range("F12").Value = 44574.86180555556

$ws.Range("B13").Value = "139"
$ws.Range("E13").Value = 44574.88263888889
$ws.Range("F13").Value = 44575.16666666666

$ws.Range("E14").Value = 44570.45763888889
$ws.Range("F14").Value = 44570.77708333333

$ws.Range("E15").Value = 44572.00625
$ws.Range("F15").Value = 44572.32569444444

$ws.Range("E16").Value = 44571.23194444444
$ws.Range("F16").Value = 44571.55138888889

$ws.Range("E17").Value = 44573.09375
$ws.Range("F17").Value = 44573.45833333334

$ws.Range("E18").Value = 44573.55347222222
$ws.Range("F18").Value = 44573.87291666667

$ws.Range("E19").Value = 44572.77986111111
$ws.Range("F19").Value = 44573.09930555556

$ws.Range("E20").Value = 44568.7875
$ws.Range("F20").Value = 44569.06875

$ws.Range("E21").Value = 44570.33611111111
$ws.Range("F21").Value = 44570.61736111111

$ws.Range("E22").Value = 44569.56180555555
$ws.Range("F22").Value = 44569.84305555555

$ws.Range("E23").Value = 44571.59241071429
$ws.Range("F23").Value = 44572.06254960317

$ws.Range("E24").Value = 44571.88402777778
$ws.Range("F24").Value = 44572.16458333333

$ws.Range("E25").Value = 44571.11041666667
$ws.Range("F25").Value = 44571.39097222222

$ws.Range("E26").Value = 44567.32152777778
$ws.Range("F26").Value = 44567.6

$ws.Range("E28").Value = 44566.66116071428
$ws.Range("F28").Value = 44566.9396329365

$ws.Range("E31").Value = 44565.33606150793
$ws.Range("F31").Value = 44565.61453373015

$ws.Range("E32").Value = 44568.30625
$ws.Range("F32").Value = 44568.58472222222

$ws.Range("E33").Value = 44569.64444444444
$ws.Range("F33").Value = 44569.92291666667

$ws.Range("E34").Value = 44568.96661706348
$ws.Range("F34").Value = 44569.24508928571

$ws.Range("E35").Value = 44565.95491071428
$ws.Range("F35").Value = 44566.24518849206

$ws.Range("E37").Value = 44562
$ws.Range("F37").Value = 44562.29027777778

$ws.Range("E40").Value = 44563.96944444445
$ws.Range("F40").Value = 44564.25972222222

$ws.Range("E41").Value = 44566.9396329365
$ws.Range("F41").Value = 44567.22991071428

$ws.Range("E42").Value = 44568.27782738095
$ws.Range("F42").Value = 44568.56810515872

$ws.Range("E43").Value = 44564.95416666667
$ws.Range("F43").Value = 44565.24444444444
